$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3  = 5
    4  = -4
    6  = -7
    15 = -2
    16 = -12
    18 = -4
    21 = -4
    23 = 1
    24 = 2
    28 = 10
    36 = -5
    38 = -2
    40 = -3
    41 = -4
    42 = -4
    43 = -3
    46 = -4
    47 = -9
    48 = -2
    50 = -4
    55 = 2
    62 = -4
    65 = -2
    66 = -4
    69 = -4
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
